$wb = $excel.ActiveWorkbook

# Update "展览" sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 30
$ws1.Range("F3").Value = 58

# Update "全部类型" sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 30
$ws4.Range("F3").Value = 58
